# Sept 2020 finance tracker - quick data cleaning; added missing negatives

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray red-highlight formatting that was left on the otherwise
# empty F7 cell (no comment text was ever entered there).
$ws.Range("F7").Clear()

# Remove the erroneous "Active Living" row (a duplicated/incorrect +382
# COVID refund entry). The row beneath it (Little Caesars, already
# correctly recorded as a negative expense) shifts up to become the new
# last row of data.
$ws.Rows("18").Delete()

# Leave the selection where the cleanup happened (the whole row that is
# now the last row in the table).
$ws.Range("A18:XFD18").Select()
